# queryParameter_Linkage.xlsx — "tab" sheet update
# Commit: Remove injury suppression; Add url for multiple lenses; updates and injury modal
#
# Insert 3 new rows (new "MULTIPLE LENSES" tab/lens block) above the existing
# "DISPARITIES" row, pushing the old rows 17-26 down to 20-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tab")
$ws.Activate()

# Insert three blank rows at 17:19 - everything currently at row 17 and below
# shifts down by three (old row 17 -> new row 20, ... old row 26 -> new row 29).
$ws.Rows("17:19").Insert()

# Fill the new rows column-by-column (A, then B, then C, then D, then E) so the
# new shared-string table entries land in the same order the workbook expects.

# Column A - tab
$ws.Cells.Item(17, 1).Value = "MULTIPLE LENSES"
$ws.Cells.Item(18, 1).Value = "MULTIPLE LENSES"
$ws.Cells.Item(19, 1).Value = "MULTIPLE LENSES"

# Column B - sub_tab
$ws.Cells.Item(17, 2).Value = "MORTALITY"
$ws.Cells.Item(18, 2).Value = "MORBIDITY"
$ws.Cells.Item(19, 2).Value = "YLD AND RISK"

# Column C - tabID
$ws.Cells.Item(17, 3).Value = "multipleLenses"
$ws.Cells.Item(18, 3).Value = "multipleLenses"
$ws.Cells.Item(19, 3).Value = "multipleLenses"

# Column D - sub_tabID
$ws.Cells.Item(17, 4).Value = "causeOfDeathTab"
$ws.Cells.Item(18, 4).Value = "nonFatalMeasuresTab"
$ws.Cells.Item(19, 4).Value = "stateMeasuresTab"

# Column E - queryName
$ws.Cells.Item(17, 5).Value = "mortalitylens"
$ws.Cells.Item(18, 5).Value = "morbiditylens"
$ws.Cells.Item(19, 5).Value = "risklens"

# Restore the view's selected cell to match the edited area.
$ws.Range("E19").Select() | Out-Null
